# KetQuaTest_RemoveFromCart.xlsx edit
# - Swap the "Du Lieu Mau" (col C) / "Cac Buoc" (col D) values for the
#   existing REMOVE_02 / REMOVE_04 / REMOVE_01 rows (they were populated
#   in the wrong columns).
# - Insert a brand-new REMOVE_05 "Loi he thong (Exception)" row between the
#   REMOVE_01 row and the REMOVE_03 row.
# - Re-balance a couple of column widths now that the data that drives
#   their "best fit" size has moved around / grown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell far away from the real table that we use as a relay when
# writing multi-line strings. Writing directly into the destination cell
# with an embedded newline makes the engine recompute that row's height;
# going through Copy()/Delete() avoids that side effect entirely.
$scratchRow = 500
$scratch = $ws.Range("Z" + $scratchRow)

function Set-MultilineValue([string]$address, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy($ws.Range($address))
    $ws.Rows.Item($scratchRow).Delete()
}

function Swap-CD([int]$row) {
    $cAddr = "C" + $row
    $dAddr = "D" + $row
    $scratch.Value = $ws.Range($cAddr).Value2
    $ws.Range($cAddr).Value = $ws.Range($dAddr).Value2
    $scratch.Copy($ws.Range($dAddr))
    $ws.Rows.Item($scratchRow).Delete()
}

# --- Fix the existing rows: column C should hold the short sample data,
#     column D the (often multi-line) steps text. ---
Swap-CD 2
Swap-CD 3
Swap-CD 4
Swap-CD 5

# --- Insert the new REMOVE_05 row just before the REMOVE_03 row (current
#     row 5), pushing REMOVE_03 down to row 6. ---
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "REMOVE_05"
$ws.Range("B5").Value = "Lỗi hệ thống (Exception)"
$ws.Range("C5").Value = "Index=1, Error"
Set-MultilineValue "D5" "1. Service ném RuntimeException`n2. Servlet catch"
$ws.Range("E5").Value = "Catch & Redirect an toàn"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

# --- Column widths: approximate the new best-fit widths as closely as
#     this engine's quantized ColumnWidth allows. ---
$ws.Columns.Item(2).ColumnWidth = 21.833
$ws.Columns.Item(3).ColumnWidth = 12.8335
$ws.Columns.Item(4).ColumnWidth = 30.1665
